$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New child rows 10 and 11: value in column G, label in column B
# (value entered first so the shared-string table matches insertion order)
$ws.Range("G10").Value = "エラーにならない～！"
$ws.Range("B10").Value = "記号！"

$ws.Range("G11").Value = "エラーにならない～～！"
$ws.Range("B11").Value = "（記／号）～"

# New "日本語！" group header (row 9), mirroring the existing "parent" / "customer" groups
$ws.Range("A9").Value = "日本語！"

# New named ranges pointing at the newly added values (wildcard matching of
# correct item names from JSON Schema)
$wb.Names.Add("json.日本語_._記_号__", $ws.Range("G11"))
$wb.Names.Add("json.日本語_.記号_", $ws.Range("G10"))

# Leave the selection on the newly added header cell
$ws.Range("A9").Select() | Out-Null
